$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '245.62'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '0.74%'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '29.06'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '5.73%'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.180'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '1.35%'

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '0.98%'

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '0.28%'

$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.108'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '3.28%'

$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8603'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '5.00%'

$ws.Range('B9').Value = 'FTXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.8615'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '0.59%'

$ws.Range('B10').Value = 'One'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0006012'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '0.08%'

$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1364'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '2.34%'

$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07064'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '1.71%'

$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03056'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '6.42%'

$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09371'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.15%'

$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001541'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '1.98%'

$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.006001'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-3.41%'

$ws.Range('B17').Value = 'UpBots'
$ws.Range('C17').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.007489'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '5,227.32%'

$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.491'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-0.63%'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.269'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '4.45%'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.3201'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '1.58%'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.03313'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '2.74%'

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.1283'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-1.46%'

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.480'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-2.89%'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.04154'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '1.33%'

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '0.52%'

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '1.04%'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.004992'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '11.66%'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0001210'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '2.53%'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03742'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '0.73%'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.005791'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '67.78%'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1071'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '1.28%'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002100'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-8.68%'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.008585'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-5.18%'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005273'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '3.20%'

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '0.06%'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05702'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-43.52%'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.002256'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '-11.22%'

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002100'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '0.06%'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0002000'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.06%'
